# Updated cryptos list on Sun Jan 28 18:58:12 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-decimal-looking Price strings to stay text (avoid Excel auto-
# converting them to numbers), matching the original inline-string cells.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D17", "D19", "D22", "D23", "D24", "D27", "D28", "D30", "D31", "D35", "D36", "D47", "D48", "D49", "D50", "D51", "D44", "D45")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.992.48"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "2.261.02"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "305.15"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").Value = "95.24"
$ws.Range("E6").Value = "  +2.56%  "

$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("D10").Value = "35.01"
$ws.Range("E10").Value = "  +7.07%  "

$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("D13").Value = "6.61"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").Value = "2.612.25"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "2.271.61"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "41.901.03"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").Value = "12.37"
$ws.Range("E19").Value = "  -4.27%  "

$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").Value = "67.61"
$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("D23").Value = "237.15"
$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("D24").Value = "2.56"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").Value = "23.66"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").Value = "36.57"
$ws.Range("E28").Value = "  +5.56%  "

$ws.Range("E29").Value = "  -1.47%  "

$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  +1.94%  "

$ws.Range("D31").Value = "160.30"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("E32").Value = "  -2.82%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  +5.10%  "

$ws.Range("D35").Value = "0.0736"
$ws.Range("E35").Value = "  -0.79%  "

$ws.Range("D36").Value = "16.99"
$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("E38").Value = "  -0.55%  "

$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("E41").Value = "  +1.44%  "

$ws.Range("E42").Value = "  +6.98%  "

$ws.Range("D43").Value = "1.979.46"
$ws.Range("E43").Value = "  -1.55%  "

# Rows 44 and 45 swap places (EnergySwap <-> VeChain)
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.78"
$ws.Range("E45").Value = "  -3.85%  "

$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("D47").Value = "9.90"
$ws.Range("E47").Value = "  -3.70%  "

$ws.Range("D48").Value = "53.09"
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("D49").Value = "72.21"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("D50").Value = "1.50"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").Value = "90.74"
$ws.Range("E51").Value = "  -1.03%  "

